# docs: close M2, start M3, and add release gate checklist
#
# - WBS sheet rows 27-49 (Phase 2 "Core DocTypes" tasks): Execution Status
#   "Partial" -> "Done", and a "Completed On" date (column K) is stamped
#   with the same date already present in "Started On" (column J).
# - WBS sheet rows 50-62 (Phase 3 "Traceability Engine" tasks): Execution
#   Status "Not Started" -> "Partial", and a "Started On" date (column J)
#   of 2026-04-17 is recorded.
# - Milestones sheet: M2 status moves to "Done", M3 status moves to
#   "In Progress".

$wb = $excel.ActiveWorkbook

$wbs = $wb.Worksheets.Item("WBS")

# Phase 2 tasks (rows 27-49): mark Done, stamp Completed On (col K) with
# the same date already recorded in Started On (col J).
for ($r = 27; $r -le 49; $r++) {
    $wbs.Cells.Item($r, 8).Value = "Done"
    $startedOn = $wbs.Cells.Item($r, 10).Value
    $wbs.Cells.Item($r, 11).Value = $startedOn
}

# Phase 3 tasks (rows 50-62): move from Not Started to Partial, and
# record the Started On date (col J).
for ($r = 50; $r -le 62; $r++) {
    $wbs.Cells.Item($r, 8).Value = "Partial"
    $wbs.Cells.Item($r, 10).Value = "2026-04-17"
}

# Milestones: M2 closed out, M3 kicked off.
$milestones = $wb.Worksheets.Item("Milestones")
$milestones.Range("F4").Value = "✅ Done"
$milestones.Range("F5").Value = "🟨 In Progress"
